# "error solve ifrs list" -- fix the IFRS financial data table for 용평리조트
# (Yongpyong Resort): the annual rows (2014-2018, rows 2-6) had the wrong
# figures (looked like they were pulled from a much larger / different
# company), so every metric column (D:AJ) on those rows is corrected to the
# right numbers. A couple of stray per-share metrics (PER/AD, 현금배당수익률/AH)
# that don't apply to rows 2-3 are cleared. The forecast rows for 2019-2021
# (rows 7-9) turn out to have been entirely bogus/placeholder data, so all of
# their metric cells (D:AJ) are wiped, leaving just the row index/label
# columns (A:C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (2014/12)
$ws.Range("D2").Value = 1336
$ws.Range("E2").Value = 172
$ws.Range("F2").Value = 172
$ws.Range("G2").Value = 54
$ws.Range("H2").Value = 40
$ws.Range("I2").Value = 40
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 8189
$ws.Range("L2").Value = 5174
$ws.Range("M2").Value = 3015
$ws.Range("N2").Value = 3015
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 1900
$ws.Range("Q2").Value = 405
$ws.Range("R2").Value = -407
$ws.Range("S2").Value = 219
$ws.Range("T2").Value = 104
$ws.Range("U2").Value = 301
$ws.Range("V2").Value = 2254
$ws.Range("W2").Value = 12.88
$ws.Range("X2").Value = 2.96
$ws.Range("Y2").Value = 1.32
$ws.Range("Z2").Value = 0.5
$ws.Range("AA2").Value = 171.61
$ws.Range("AB2").Value = 59.15
$ws.Range("AC2").Value = 104
$ws.Range("AD2").ClearContents()
$ws.Range("AE2").Value = 7933
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = 0
$ws.Range("AH2").ClearContents()
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 38000000

# Row 3 (2015/12)
$ws.Range("D3").Value = 1763
$ws.Range("E3").Value = 264
$ws.Range("F3").Value = 264
$ws.Range("G3").Value = 152
$ws.Range("H3").Value = 116
$ws.Range("I3").Value = 116
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 9707
$ws.Range("L3").Value = 6578
$ws.Range("M3").Value = 3128
$ws.Range("N3").Value = 3128
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 1900
$ws.Range("Q3").Value = 368
$ws.Range("R3").Value = -71
$ws.Range("S3").Value = -309
$ws.Range("T3").Value = 136
$ws.Range("U3").Value = 232
$ws.Range("V3").Value = 1703
$ws.Range("W3").Value = 14.99
$ws.Range("X3").Value = 6.56
$ws.Range("Y3").Value = 3.78
$ws.Range("Z3").Value = 1.29
$ws.Range("AA3").Value = 210.29
$ws.Range("AB3").Value = 65.11
$ws.Range("AC3").Value = 306
$ws.Range("AD3").ClearContents()
$ws.Range("AE3").Value = 8231
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = 0
$ws.Range("AH3").ClearContents()
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 38000000

# Row 4 (2016/12)
$ws.Range("D4").Value = 2107
$ws.Range("E4").Value = 304
$ws.Range("F4").Value = 304
$ws.Range("G4").Value = 182
$ws.Range("H4").Value = 143
$ws.Range("I4").Value = 143
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 9289
$ws.Range("L4").Value = 5325
$ws.Range("M4").Value = 3964
$ws.Range("N4").Value = 3964
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 2407
$ws.Range("Q4").Value = -190
$ws.Range("R4").Value = 156
$ws.Range("S4").Value = 296
$ws.Range("T4").Value = 178
$ws.Range("U4").Value = -368
$ws.Range("V4").Value = 1300
$ws.Range("W4").Value = 14.42
$ws.Range("X4").Value = 6.8
$ws.Range("Y4").Value = 4.04
$ws.Range("Z4").Value = 1.51
$ws.Range("AA4").Value = 134.32
$ws.Range("AB4").Value = 64.72
$ws.Range("AC4").Value = 324
$ws.Range("AD4").Value = 26.52
$ws.Range("AE4").Value = 8236
$ws.Range("AF4").Value = 1.04
$ws.Range("AG4").Value = 100
$ws.Range("AH4").Value = 1.17
$ws.Range("AI4").Value = 33.61
$ws.Range("AJ4").Value = 48133333

# Row 5 (2017/12)
$ws.Range("D5").Value = 2097
$ws.Range("E5").Value = 335
$ws.Range("F5").Value = 335
$ws.Range("G5").Value = 219
$ws.Range("H5").Value = 170
$ws.Range("I5").Value = 170
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 8682
$ws.Range("L5").Value = 4600
$ws.Range("M5").Value = 4082
$ws.Range("N5").Value = 4082
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 2407
$ws.Range("Q5").Value = 178
$ws.Range("R5").Value = -615
$ws.Range("S5").Value = 17
$ws.Range("T5").Value = 584
$ws.Range("U5").Value = -406
$ws.Range("V5").Value = 1424
$ws.Range("W5").Value = 15.97
$ws.Range("X5").Value = 8.09
$ws.Range("Y5").Value = 4.21
$ws.Range("Z5").Value = 1.89
$ws.Range("AA5").Value = 112.69
$ws.Range("AB5").Value = 69.6
$ws.Range("AC5").Value = 352
$ws.Range("AD5").Value = 29.67
$ws.Range("AE5").Value = 8480
$ws.Range("AF5").Value = 1.23
$ws.Range("AG5").Value = 100
$ws.Range("AH5").Value = 0.96
$ws.Range("AI5").Value = 28.39
$ws.Range("AJ5").Value = 48133333

# Row 6 (2018/12)
$ws.Range("D6").Value = 1814
$ws.Range("E6").Value = 246
$ws.Range("F6").Value = 246
$ws.Range("G6").Value = 131
$ws.Range("H6").Value = 102
$ws.Range("I6").Value = 102
$ws.Range("K6").Value = 8128
$ws.Range("L6").Value = 3986
$ws.Range("M6").Value = 4142
$ws.Range("N6").Value = 4142
$ws.Range("P6").Value = 2407
$ws.Range("Q6").Value = 525
$ws.Range("R6").Value = -172
$ws.Range("S6").Value = -81
$ws.Range("T6").Value = 162
$ws.Range("U6").Value = 363
$ws.Range("V6").Value = 1546
$ws.Range("W6").Value = 13.54
$ws.Range("X6").Value = 5.62
$ws.Range("Y6").Value = 2.48
$ws.Range("Z6").Value = 1.21
$ws.Range("AA6").Value = 96.24
$ws.Range("AB6").Value = 71.72
$ws.Range("AC6").Value = 212
$ws.Range("AD6").Value = 28.77
$ws.Range("AE6").Value = 8604
$ws.Range("AF6").Value = 0.71
$ws.Range("AG6").Value = 70
$ws.Range("AH6").Value = 1.15
$ws.Range("AI6").Value = 33.07
$ws.Range("AJ6").Value = 48133333

# Rows 7-9 (2019E/2020E/2021E): the forecast figures were placeholder/bogus
# data -- clear all metric columns (D:AJ), keep the A/B/C label columns.
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
